# "Changed the launch intent button to be Record Location"
#
# The "survey" sheet is an XLSForm-style sheet. Add a new "buttonLabel"
# column (F) with the header "buttonLabel" and give the "Location"
# geopoint row (row 5) a button label of "Record Location" so the
# launch-intent button in the generated form shows that text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")

# New column header + value for the geopoint row's button label.
$ws1.Range("F1").Value = "buttonLabel"
$ws1.Range("F5").Value = "Record Location"

# Give the new column a sensible custom width (matches the authored width
# of 14.5 "characters" once Excel's column-width padding is accounted for).
$ws1.Columns.Item(6).ColumnWidth = 13.666666666666666

# The "survey" sheet becomes the active tab/selection (it was "settings"
# before), with the cursor sitting just past the new data in column F.
$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
